$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new text value. Values are assigned with a leading
# apostrophe so Excel stores them as text (matching the source data,
# which holds price/percentage figures as text, not numbers).
$ws.Range("D2").Value = "'290.37"
$ws.Range("E2").Value = "'-3.55%"
$ws.Range("E3").Value = "'-4.91%"
$ws.Range("E4").Value = "'-0.03%"
$ws.Range("D5").Value = "'0.07212"
$ws.Range("E5").Value = "'-5.60%"
$ws.Range("D6").Value = "'1.837"
$ws.Range("E6").Value = "'-5.43%"
$ws.Range("D7").Value = "'7.691"
$ws.Range("E7").Value = "'-1.93%"
$ws.Range("D8").Value = "'3.767"
$ws.Range("E8").Value = "'-0.82%"
$ws.Range("D9").Value = "'0.8974"
$ws.Range("E9").Value = "'-2.12%"
$ws.Range("D10").Value = "'0.1662"
$ws.Range("E10").Value = "'-4.99%"
$ws.Range("D11").Value = "'0.07734"
$ws.Range("E11").Value = "'-0.82%"
$ws.Range("D12").Value = "'0.08015"
$ws.Range("E12").Value = "'-5.61%"
$ws.Range("D13").Value = "'0.03038"
$ws.Range("E13").Value = "'-3.97%"
$ws.Range("E14").Value = "'0.19%"
$ws.Range("D15").Value = "'0.001491"
$ws.Range("E15").Value = "'-1.25%"
$ws.Range("D16").Value = "'0.005857"
$ws.Range("E16").Value = "'-1.09%"
$ws.Range("D18").Value = "'3.466"
$ws.Range("E18").Value = "'0.13%"
$ws.Range("D19").Value = "'2.080"
$ws.Range("E19").Value = "'-3.30%"
$ws.Range("E20").Value = "'-0.83%"
$ws.Range("E21").Value = "'-1.61%"
$ws.Range("D22").Value = "'4.051"
$ws.Range("E22").Value = "'-5.00%"
$ws.Range("D23").Value = "'0.2391"
$ws.Range("E23").Value = "'20.08%"
$ws.Range("D24").Value = "'0.04514"
$ws.Range("E24").Value = "'-0.08%"
$ws.Range("E25").Value = "'-0.49%"
$ws.Range("D26").Value = "'0.004625"
$ws.Range("E26").Value = "'5.36%"
$ws.Range("E27").Value = "'0.11%"
$ws.Range("D39").Value = "'0.01568"
$ws.Range("E39").Value = "'-8.00%"
$ws.Range("D40").Value = "'0.04383"
$ws.Range("E40").Value = "'-6.24%"
$ws.Range("D41").Value = "'0.007318"
$ws.Range("E41").Value = "'-2.09%"
$ws.Range("D42").Value = "'0.009906"
$ws.Range("D43").Value = "'0.1302"
$ws.Range("E43").Value = "'-3.44%"
$ws.Range("D44").Value = "'0.002016"
$ws.Range("E44").Value = "'-13.52%"
$ws.Range("D45").Value = "'0.009531"
$ws.Range("E45").Value = "'-8.86%"
$ws.Range("D46").Value = "'0.00005963"
$ws.Range("E46").Value = "'-4.78%"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'0.11%"
$ws.Range("E48").Value = "'174.74%"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("E50").Value = "'0.11%"
$ws.Range("D51").Value = "'0.0002003"
$ws.Range("E51").Value = "'0.11%"
